# "versao aperfeicoada para 2 tipos de emissao de boletos"
# Expand the one-client sheet into a multi-client billing list, storing the
# ID column as text (so long numeric IDs keep their leading digits/format)
# while re-using the existing "ID" header look for the data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows that go under the existing header (Cliente / ID / Valor).
$data = @(
    @("maxi massas", "704083103754311", 25),
    @("Casa Deliza AQA/SC", "4766441153460561", 1000),
    @("Barão Consórcios", "332223135566070", 400),
    @("Shopping Lupo", "866541220063713", 400),
    @("Borsari Imóveis", "368957994998298", 400),
    @("Agrotécnica", "538295984311594", 500),
    @("Animalia Rio Preto", "690102435386578", 400),
    @("Passarinho Hortifrúti", "957297151620080", 500),
    @("Julianeti", "725054975582041", 1000),
    @("Casa9", "1314936005979662", 1100),
    @("Spazzeo", "350829366163255", 740),
    @("Hotel Salto Grande", "757419351915951", 400),
    @("Fuba", "747162889648020", 400),
    @("Dental Power", "1119662528586518", 600),
    @("Apoio", "477908569850788", 1000),
    @("Micelli", "618129976369036", 200),
    @("Trinity", "589096946057116", 300)
)

$lastRow = $data.Length + 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]

    if ($r -eq $lastRow) {
        # Last row's ID cell was formatted by hand instead of being copied
        # from the header, so it ends up on its own (near-identical) style.
        $ws.Cells.Item($r, 2).NumberFormat = "@"
        $ws.Cells.Item($r, 2).Font.Name = "JetBrains Mono"
        $ws.Cells.Item($r, 2).Font.Size = 10
        $ws.Cells.Item($r, 2).Font.Color = 8421504
        $ws.Cells.Item($r, 2).VerticalAlignment = -4108
        $ws.Cells.Item($r, 2).Value = $row[1]
    } else {
        # Re-use the header's "ID" look (JetBrains Mono, grey,
        # vertical-centered) for every ID cell, then force text storage so
        # the long numbers are kept verbatim instead of being parsed as
        # numbers.
        $ws.Range("B1").Copy()
        $ws.Cells.Item($r, 2).PasteSpecial(-4122)
        $ws.Cells.Item($r, 2).NumberFormat = "@"
        $ws.Cells.Item($r, 2).Value = $row[1]
    }

    $ws.Cells.Item($r, 3).Value = $row[2]
}

# A stray formatting touch on D2 (underline) left over from the edit.
$ws.Range("D2").Font.Underline = 2

# Match columns A and B widths (client name / id) to fit the new content.
$ws.Range("A1").ColumnWidth = 18.6
$ws.Range("B1").ColumnWidth = 18.6

[void]$ws.Range("A2:C3").Select()
